$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before "总计", with the same layout
#    (header row + column-A numbering style) as the "2021-Q4" sheet, so the
#    new sheet inherits identical formatting/styles.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$q4.Copy($totalSheet)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# The duplicated sheet only has 4 data rows (rows 2-5 incl. header); we need
# 7 data rows (rows 2-8). Insert 4 more rows and copy the formatting of the
# last existing data row down onto them.
$new.Rows("5:8").Insert()
$new.Cells.Item(4, 1).Copy($new.Cells.Item(5, 1))
$new.Cells.Item(4, 1).Copy($new.Cells.Item(6, 1))
$new.Cells.Item(4, 1).Copy($new.Cells.Item(7, 1))
$new.Cells.Item(4, 1).Copy($new.Cells.Item(8, 1))

# Columns B-G hold text values (fund codes with leading zeros, decimal
# strings, etc.) - force text format while writing so they aren't silently
# coerced to numbers, then clear the temporary formatting again.
$textRange = $new.Range("B2:G8")
$textRange.NumberFormat = "@"

$fundData = @(
    @(0, "010864", "泓德卓远混合A",       "39.08", "91.88", "2.76", "1.0786", 10),
    @(1, "008809", "安信民稳增长混合A",     "15.09", "44.77", "3.89", "0.5870", 3),
    @(2, "012256", "安信丰穗一年持有混合A", "26.49", "20.42", "1.72", "0.4556", 3),
    @(3, "009849", "安信稳健聚申一年持有期混合", "12.33", "30.81", "2.98", "0.3674", 3),
    @(4, "010865", "泓德卓远混合C",       "12.25", "91.88", "2.76", "0.3381", 10),
    @(5, "008810", "安信民稳增长混合C",     "6.76", "44.77", "3.89", "0.2630", 3),
    @(6, "012257", "安信丰穗一年持有混合C", "2.41", "20.42", "1.72", "0.0415", 3)
)

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = $i + 2
    $rec = $fundData[$i]
    $new.Cells.Item($row, 1).Value = $rec[0]
    $new.Cells.Item($row, 2).Value = $rec[1]
    $new.Cells.Item($row, 3).Value = $rec[2]
    $new.Cells.Item($row, 4).Value = $rec[3]
    $new.Cells.Item($row, 5).Value = $rec[4]
    $new.Cells.Item($row, 6).Value = $rec[5]
    $new.Cells.Item($row, 7).Value = $rec[6]
    $new.Cells.Item($row, 8).Value = $rec[7]
}

$textRange.ClearFormats()

# ---------------------------------------------------------------------------
# 2) Add a new first data row to "总计" summarising the 2022-Q1 sheet
#    (7 funds held, 3.13 亿元 total market value), shifting the existing
#    rows down by one.
# Re-fetch the "总计" worksheet by name: the sheet-collection positions
# shifted once the new sheet was inserted in front of it, so any reference
# captured earlier would now point at the wrong sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(3, 1).Copy($totalSheet.Cells.Item(2, 1))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 7
$totalSheet.Cells.Item(2, 4).Value = 3.13
